# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first data sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 633
$ws1.Range("F3").Value  = 2209
$ws1.Range("F4").Value  = 90
$ws1.Range("F5").Value  = 13204
$ws1.Range("F9").Value  = 483
$ws1.Range("F10").Value = 1187
$ws1.Range("F12").Value = 13786
$ws1.Range("F13").Value = 14403
$ws1.Range("F15").Value = 171
$ws1.Range("F21").Value = 38
$ws1.Range("F22").Value = 1095
$ws1.Range("F25").Value = 5453
$ws1.Range("F27").Value = 88
$ws1.Range("F28").Value = 331
$ws1.Range("F29").Value = 24
$ws1.Range("F30").Value = 68

# --- Sheet "全部类型" (aggregate sheet with an extra row vs. "展览") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 633
$ws4.Range("F3").Value  = 2209
$ws4.Range("F4").Value  = 90
$ws4.Range("F5").Value  = 13204
$ws4.Range("F10").Value = 483
$ws4.Range("F11").Value = 1187
$ws4.Range("F13").Value = 13786
$ws4.Range("F14").Value = 14403
$ws4.Range("F16").Value = 171
$ws4.Range("F22").Value = 38
$ws4.Range("F23").Value = 1095
$ws4.Range("F26").Value = 5454
$ws4.Range("F28").Value = 88
$ws4.Range("F29").Value = 331
$ws4.Range("F30").Value = 24
$ws4.Range("F31").Value = 68
